$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) "Formula: Paired Samples" block - replace the centered oMathPara
#    (d = D-bar / s_D) with a plain-text source-code style line.
# ------------------------------------------------------------------
$formulaPara = $d.Paragraphs(6)
$formulaRange = $formulaPara.Range

$formulaXml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="SourceCode"/></w:pPr><w:r><w:t xml:space="preserve">d = mean(differences) / std(differences)</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$formulaRange.InsertXML($formulaXml)

# Re-fetch the paragraph (content changed) and apply the VerbatimChar
# run style to the run's text (exclude the trailing paragraph mark).
$formulaPara = $d.Paragraphs(6)
$codeRange = $d.Range($formulaPara.Range.Start, $formulaPara.Range.End - 1)
$codeRange.Style = "VerbatimChar"

# ------------------------------------------------------------------
# 2) "Where:" block - replace the two inline oMath runs (D-bar, s_D)
#    with bold plain-text labels and revise the explanatory wording.
# ------------------------------------------------------------------
$wherePara = $d.Paragraphs(7)
$whereRange = $wherePara.Range

$whereXml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="FirstParagraph"/></w:pPr><w:r><w:t xml:space="preserve">Where:</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">-</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve">differences</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">= Model A F1 score minus Model B F1 score, computed for each field</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">-</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve">mean(differences)</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">= Average of these difference scores across all fields</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">-</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve">std(differences)</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">= Standard deviation of the difference scores</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$whereRange.InsertXML($whereXml)
